$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename header labels: "<name>_old" -> "<name>_FV2310" (columns A-J),
#    "<name>_new" -> "<name>_FV2404" (columns L-U). Column K ("diff") is
#    unchanged.
# ---------------------------------------------------------------------------
$baseHeaders = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$leftCols  = @("A","B","C","D","E","F","G","H","I","J")
$rightCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $baseHeaders.Length; $i++) {
    $ws.Range($leftCols[$i]  + "1").Value = $baseHeaders[$i] + "_FV2310"
    $ws.Range($rightCols[$i] + "1").Value = $baseHeaders[$i] + "_FV2404"
}

# ---------------------------------------------------------------------------
# 2) Turn the used range into an Excel Table ("Table1") without letting the
#    header row's pre-existing bold/fill formatting get captured as a new
#    dxf (that would otherwise add a headerRowDxfId + grow styles.xml).
#    Trick: stash the header formatting on a scratch row, strip formats from
#    the header row, add the table, then restore the formatting and wipe the
#    scratch row again.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A100:U100")

$scratchRange.Value = "x"
$headerRange.Copy()
$scratchRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$headerRange.ClearFormats()

$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U56"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)

$scratchRange.Copy()
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$scratchRange.Clear()

try { $lo.TableStyle = "" } catch { }

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1, top-left cell of the
#    scrolling pane is A2).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
